$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append three new key/translation rows after the existing data (rows 75-77)
$ws.Cells.Item(75, 1).Value = "LOGINSUCESSFUL"
$ws.Cells.Item(75, 2).Value = "Zalogowano pomyślnie!"
$ws.Cells.Item(75, 3).Value = "Login Sucessful!"

$ws.Cells.Item(76, 1).Value = "PASSWORDTOOSHORT"
$ws.Cells.Item(76, 2).Value = "Hasło musi się składać co najmniej z 6 znaków"
$ws.Cells.Item(76, 3).Value = "The password must consist of at least 6 characters"

$ws.Cells.Item(77, 1).Value = "ACCOUNTNOTFOUND"
$ws.Cells.Item(77, 2).Value = "Takie konto nie istnieje"
$ws.Cells.Item(77, 3).Value = "Account does not exist"

# Match the bold "Key" column style used by column A for the rest of the table
$ws.Range("A75:A77").Font.Bold = $true

# Update the view: scroll/selection state recorded in the saved file
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 61
$ws.Range("C78").Select()
